$wb = $excel.ActiveWorkbook
$excel.Calculate()
$ws = $wb.Worksheets.Item("values")
$r1 = $ws.Range("B5")
Write-Host $r1.Value()
$r2 = $ws.Range("B7")
Write-Host $r2.Value()
